$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column N (VOTO), shifting VOTO and
# RESPONSABLE 1X10 one column to the right, to make room for the new
# "TELEFONO" (phone number) column requested for the Elector data group.
$ws.Columns("N").Insert()

# Header text for the newly inserted column.
$ws.Range("N4").Value = "TELÉFONO"

# Match the formatting (font/border/alignment) of the sibling data columns
# (H:J, M) for the new column's header/body cells.
$ws.Range("H1").Copy()
$ws.Range("N1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

$ws.Range("M4").Copy()
$ws.Range("N4").PasteSpecial(-4122)

$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)

# Re-set the header text (PasteSpecial of formats only keeps value, but be
# explicit/safe in case PasteSpecial(-4122) ever touched it).
$ws.Range("N4").Value = "TELÉFONO"

# Give the new column the same width used by the other detail columns.
$ws.Columns("N").ColumnWidth = 19.166666666666668

# Extend the "DATOS DE PEQUIVEN (COMPLEJO Y GERENCIA)" merged header so it
# now spans the CÓDIGO/NOMBRE/TELÉFONO/VOTO columns (L:O).
$ws.Range("L3:O3").Merge()

# Re-apply the AutoFilter over the widened header row (A4:P4). Calling
# AutoFilter on a range that is already inside an AutoFilter toggles it off,
# so make sure it is off first.
$ws.AutoFilterMode = $false
$ws.Range("A4:P4").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Cuadro!_FilterDatabase") {
        $n.RefersTo = "=Cuadro!`$A`$4:`$P`$4"
    }
}

# Put the selection where a user would naturally land after typing the new
# header.
$ws.Range("N2").Select()
